$wb = $excel.ActiveWorkbook

# --- Sheet 1: Translations ---
$ws1 = $wb.Worksheets.Item("Translations")

# Remember column A's width so the new "Variable" column can match it.
$ws1ColAWidth = $ws1.Columns.Item(1).ColumnWidth

# Insert a new column before column B ("Type") to hold the new "Variable" column
$ws1.Columns.Item(2).Insert()

$ws1.Range("B1").Value = "Variable"
$ws1.Range("B2").Value = "s1"
$ws1.Range("B3").Value = "s1"
$ws1.Range("B4").Value = "s1"
$ws1.Range("B5").Value = "s1"

$ws1.Columns.Item(2).ColumnWidth = $ws1ColAWidth

# --- Sheet 2: @@_question ---
$ws2 = $wb.Worksheets.Item("@@_question")

$ws2.Columns.Item(2).Insert()

$ws2.Range("B1").Value = "Variable"
$ws2.Range("B2").Value = "s1"
$ws2.Range("B2").NumberFormat = "@"

# Size the new column to fit its (short) content, like the rest of the sheet's columns.
$ws2.Columns.Item(2).AutoFit()

$ws2.Activate() | Out-Null
$ws2.Range("B3").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("B6").Select() | Out-Null
